# Weekly refresh of the "Bruselas (repollito)" price series:
# a new observation (D=44838, the most recent date) is inserted as row 25,
# pushing every existing record (former rows 25-73) down by one row, so the
# used range grows from A1:R73 to A1:R74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 25..73 down to 26..74, carrying formats/styles along.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new observation.
$ws.Cells.Item(25, 1).Value  = 9
$ws.Cells.Item(25, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value  = "Metropolitana"
$ws.Cells.Item(25, 4).Value  = 44838
$ws.Cells.Item(25, 5).Value  = 13
$ws.Cells.Item(25, 6).Value  = 100112035
$ws.Cells.Item(25, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(25, 8).Value  = "Sin especificar"
$ws.Cells.Item(25, 9).Value  = "Primera"
$ws.Cells.Item(25, 10).Value = 52
$ws.Cells.Item(25, 11).Value = 17000
$ws.Cells.Item(25, 12).Value = 17000
$ws.Cells.Item(25, 13).Value = 17000
$ws.Cells.Item(25, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(25, 15).Value = "Hijuelas"
$ws.Cells.Item(25, 16).Value = 1133
$ws.Cells.Item(25, 17).Value = 15
$ws.Cells.Item(25, 18).Value = "Hortaliza"
